# Marine protists workbook update:
#  - add two new worksheets ("Turley, P. J. Mackie " and "Herndl") after
#    the existing "de Vargas" sheet, populate them with data from
#    Turley/Mackie and Herndl, and make "Herndl" the active sheet.
#  - harmonise the text style used in column B (rows 3-13) of "de Vargas"
#    with the rest of that column.

$wb = $excel.ActiveWorkbook
$deVargas = $wb.Worksheets.Item(1)

function Copy-PageSetup($srcSheet, $dstSheet) {
    $dstSheet.PageSetup.LeftMargin = $srcSheet.PageSetup.LeftMargin
    $dstSheet.PageSetup.RightMargin = $srcSheet.PageSetup.RightMargin
    $dstSheet.PageSetup.TopMargin = $srcSheet.PageSetup.TopMargin
    $dstSheet.PageSetup.BottomMargin = $srcSheet.PageSetup.BottomMargin
    $dstSheet.PageSetup.HeaderMargin = $srcSheet.PageSetup.HeaderMargin
    $dstSheet.PageSetup.FooterMargin = $srcSheet.PageSetup.FooterMargin
    $dstSheet.PageSetup.PaperSize = $srcSheet.PageSetup.PaperSize
    $dstSheet.PageSetup.Orientation = $srcSheet.PageSetup.Orientation
    $dstSheet.PageSetup.Zoom = $srcSheet.PageSetup.Zoom
    $dstSheet.PageSetup.CenterHeader = $srcSheet.PageSetup.CenterHeader
    $dstSheet.PageSetup.CenterFooter = $srcSheet.PageSetup.CenterFooter
}

# ---------------------------------------------------------------------
# 1. "de Vargas" - small style cleanup on B3:B13 (align with B14:B.. style)
# ---------------------------------------------------------------------
for ($r = 3; $r -le 13; $r++) {
    $cell = $deVargas.Cells.Item($r, 2)
    $cell.NumberFormat = "General"
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
}

# ---------------------------------------------------------------------
# 2. Add "Turley, P. J. Mackie " worksheet right after "de Vargas"
# ---------------------------------------------------------------------
$turley = $wb.Worksheets.Add($null, $deVargas)
$turley.Name = "Turley, P. J. Mackie "

$turley.Cells.Item(2, 1).Value = "Depth [m]"
$turley.Cells.Item(2, 2).Value = "Concentration of Bacteria [cells mL^-1]"
$turley.Cells.Item(2, 3).Value = "Concentration of Flagellates [cells mL^-1]"

$turleyRows = @(
    @(45, 210000000, $null),
    @(45, 1730000000, 33000000),
    @(55, 2540000000, 29900000),
    @(45, 820000000, 3200000),
    @(45, 250000000, 1300000),
    @(300, 30000000, $null),
    @(45, 570000000, $null),
    @(300, 2110000000, $null),
    @(45, 610000000, $null),
    @(300, 300000000, 2900000),
    @(50, 340000000, 3800000),
    @(100, 130000000, 1000000),
    @(200, 510000000, 4500000),
    @(300, 130000000, 1500000),
    @(50, 3780000000, 4000000),
    @(100, 2867000000, 4300000),
    @(200, 3390000000, 2500000),
    @(300, 2810000000, 1000000),
    @(45, 3760000000, 200000)
)

$r = 3
foreach ($row in $turleyRows) {
    $turley.Cells.Item($r, 1).Value = $row[0]

    $bCell = $turley.Cells.Item($r, 2)
    $bCell.Value = $row[1]
    $bCell.NumberFormat = "0.00E+00"

    if ($null -ne $row[2]) {
        $cCell = $turley.Cells.Item($r, 3)
        $cCell.Value = $row[2]
        $cCell.NumberFormat = "0.00E+00"
    }

    $r = $r + 1
}

Copy-PageSetup $deVargas $turley
$turley.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Add "Herndl" worksheet right after "Turley, P. J. Mackie "
# ---------------------------------------------------------------------
$herndl = $wb.Worksheets.Add($null, $turley)
$herndl.Name = "Herndl"

$herndl.Cells.Item(1, 1).Value = "Reference"
$herndl.Cells.Item(1, 2).Value = "Values taken from Table 1 in Herndl. Carbon concentrations were derived by using measured cell volumes and conversion factor described in the paper."

$herndl.Cells.Item(2, 1).Value = "Concentration of Bacteria [cells mL^-1]"
$herndl.Cells.Item(2, 2).Value = "Concentration of Flagellates [cells mL^-1]"
$herndl.Cells.Item(2, 3).Value = "Concentration of Monads [cells mL^-1]"
$herndl.Cells.Item(2, 4).Value = "Bacterial carbon concentration [g C mL^-1]"
$herndl.Cells.Item(2, 5).Value = "Flagellates carbon concentration [g C mL^-1]"
$herndl.Cells.Item(2, 6).Value = "Monads carbon concentration [g C mL^-1]"

$herndlRows = @(
    @(55900000, 1750000, 400000),
    @(15200000, 1290000, 840000),
    @(2500000, 30000, 30000),
    @(5800000, 180000, 40000),
    @(19500000, 280000, 80000),
    @(1300000, 20000, 20000),
    @(110400000, 2070000, 1160000),
    @(2100000, 30000, 20000),
    @(35600000, 520000, 330000)
)

$r = 3
foreach ($row in $herndlRows) {
    $aCell = $herndl.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.NumberFormat = "0.00E+00"

    $bCell = $herndl.Cells.Item($r, 2)
    $bCell.Value = $row[1]
    $bCell.NumberFormat = "0.00E+00"

    $cCell = $herndl.Cells.Item($r, 3)
    $cCell.Value = $row[2]
    $cCell.NumberFormat = "0.00E+00"

    $dCell = $herndl.Cells.Item($r, 4)
    $dCell.Formula = "=AVERAGE(0.067,0.25)*380*A$r*1E-015"
    $dCell.NumberFormat = "General"

    $eCell = $herndl.Cells.Item($r, 5)
    $eCell.Formula = "=11.43*220*B$r*1E-015"
    $eCell.NumberFormat = "General"

    $fCell = $herndl.Cells.Item($r, 6)
    $fCell.Formula = "=2.46*220*C$r*1E-015"
    $fCell.NumberFormat = "General"

    $r = $r + 1
}

Copy-PageSetup $deVargas $herndl
$herndl.Range("A13").Select()

# ---------------------------------------------------------------------
# 4. "Herndl" becomes the active tab of the workbook
# ---------------------------------------------------------------------
$herndl.Activate()
